# Auto-generated Excel COM-interop script to apply cryptos.xlsx data refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($ws, [string]$cellRef, [string]$val)
    $rng = $ws.Range($cellRef)
    # Force a Text number format before assigning so Excel does not
    # reinterpret numeric-looking strings (e.g. '1.00', '3.20') as numbers,
    # then clear the format again so no new cell style is introduced.
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.ClearFormats()
}

Set-TextCell $ws 'D2' '61.616.00'
Set-TextCell $ws 'E2' '  -2.95%  '
Set-TextCell $ws 'E3' '  -2.88%  '
Set-TextCell $ws 'E4' '  +0.04%  '
Set-TextCell $ws 'D5' '405.33'
Set-TextCell $ws 'E5' '  -2.81%  '
Set-TextCell $ws 'D6' '133.63'
Set-TextCell $ws 'E6' '  +8.68%  '
Set-TextCell $ws 'E7' '  -2.33%  '
Set-TextCell $ws 'E8' '  +0.05%  '
Set-TextCell $ws 'D9' '0.669'
Set-TextCell $ws 'E9' '  -2.57%  '
Set-TextCell $ws 'E10' '  -8.81%  '
Set-TextCell $ws 'D11' '42.45'
Set-TextCell $ws 'E11' '  +2.62%  '
Set-TextCell $ws 'E12' '  -1.08%  '
Set-TextCell $ws 'D13' '3.905.18'
Set-TextCell $ws 'E13' '  -3.26%  '
Set-TextCell $ws 'D14' '8.40'
Set-TextCell $ws 'E14' '  -2.53%  '
Set-TextCell $ws 'D15' '19.76'
Set-TextCell $ws 'E15' '  -1.28%  '
Set-TextCell $ws 'D16' '3.386.99'
Set-TextCell $ws 'E16' '  -2.77%  '
Set-TextCell $ws 'D17' '61.536.62'
Set-TextCell $ws 'E17' '  -2.66%  '
Set-TextCell $ws 'E18' '  -1.65%  '
Set-TextCell $ws 'D19' '11.09'
Set-TextCell $ws 'E19' '  +0.43%  '
Set-TextCell $ws 'D20' '0.0000128'
Set-TextCell $ws 'E20' '  -11.44%  '
Set-TextCell $ws 'D21' '3.20'
Set-TextCell $ws 'E21' '  -3.61%  '
Set-TextCell $ws 'D22' '85.42'
Set-TextCell $ws 'E22' '  +2.76%  '
Set-TextCell $ws 'D23' '316.31'
Set-TextCell $ws 'E23' '  +0.14%  '
Set-TextCell $ws 'D24' '12.78'
Set-TextCell $ws 'E24' '  -0.89%  '
Set-TextCell $ws 'E25' '  -1.49%  '
Set-TextCell $ws 'D26' '4.79'
Set-TextCell $ws 'E26' '  +11.00%  '
Set-TextCell $ws 'D27' '29.51'
Set-TextCell $ws 'E27' '  -5.76%  '
Set-TextCell $ws 'D28' '8.31'
Set-TextCell $ws 'E28' '  +5.16%  '
Set-TextCell $ws 'D29' '7.58'
Set-TextCell $ws 'E29' '  -1.98%  '
Set-TextCell $ws 'B30' 'Hedera'
Set-TextCell $ws 'C30' 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextCell $ws 'D30' '0.116'
Set-TextCell $ws 'E30' '  -0.07%  '
Set-TextCell $ws 'B31' 'Kaspa'
Set-TextCell $ws 'C31' 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextCell $ws 'D31' '0.171'
Set-TextCell $ws 'E31' '  -2.94%  '
Set-TextCell $ws 'B32' 'Toncoin'
Set-TextCell $ws 'C32' 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-TextCell $ws 'D32' '2.65'
Set-TextCell $ws 'E32' '  +2.85%  '
Set-TextCell $ws 'D33' '11.35'
Set-TextCell $ws 'E33' '  -2.50%  '
Set-TextCell $ws 'D34' '1.00'
Set-TextCell $ws 'E34' '  -0.67%  '
Set-TextCell $ws 'D35' '41.08'
Set-TextCell $ws 'E35' '  -2.15%  '
Set-TextCell $ws 'D36' '0.0478'
Set-TextCell $ws 'E36' '  -2.31%  '
Set-TextCell $ws 'D37' '51.64'
Set-TextCell $ws 'E37' '  -1.32%  '
Set-TextCell $ws 'D38' '0.998'
Set-TextCell $ws 'E38' '  +0.20%  '
Set-TextCell $ws 'D39' '3.41'
Set-TextCell $ws 'E39' '  -2.19%  '
Set-TextCell $ws 'E40' '  -3.77%  '
Set-TextCell $ws 'D41' '139.85'
Set-TextCell $ws 'E41' '  +3.13%  '
Set-TextCell $ws 'D42' '1.98'
Set-TextCell $ws 'E42' '  -1.73%  '
Set-TextCell $ws 'E43' '  -2.10%  '
Set-TextCell $ws 'E44' '  +4.35%  '
Set-TextCell $ws 'D45' '3.99'
Set-TextCell $ws 'E45' '  +1.82%  '
Set-TextCell $ws 'D46' '16.60'
Set-TextCell $ws 'E46' '  -2.59%  '
Set-TextCell $ws 'E47' '  -1.43%  '
Set-TextCell $ws 'D48' '21.19'
Set-TextCell $ws 'E48' '  -3.91%  '
Set-TextCell $ws 'D49' '2.117.24'
Set-TextCell $ws 'E49' '  -3.22%  '
Set-TextCell $ws 'D50' '2.29'
Set-TextCell $ws 'E50' '  -5.82%  '
Set-TextCell $ws 'E51' '  -0.33%  '
